$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the latest crypto data refresh.
# D-column values are forced through a Text number format while assigning so that
# numeric-looking strings (e.g. "1.002", "0.3900") are stored as literal text rather
# than being auto-converted/truncated by Excel, then the format is restored to General.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.976.84'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = '  -1.98%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.651.56'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = '  -0.92%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = '  +0.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.89'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -1.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3900'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  -1.24%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3817'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  -2.45%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '52.25'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  +0.85%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.349'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  -4.13%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.002'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +0.32%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08445'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  -1.60%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.83'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  -2.54%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.067'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  -3.52%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.025'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  +1.75%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001308'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  -2.46%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.647.10'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  -0.89%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.51'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  -1.11%  '

$ws.Range("E19").Value = '  +0.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.69'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  -4.14%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.973'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  -0.37%  '

$ws.Range("E22").Value = '  +0.22%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.79'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +0.56%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.968.55'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  -2.05%  '

$ws.Range("E25").Value = '  +0.34%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.979'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  -1.19%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.10'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  -1.91%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '152.61'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  -3.33%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.420'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  +0.97%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '138.01'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -3.32%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.974'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  -1.34%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.523'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  +0.42%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.824.06'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  -0.86%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.025'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  -3.74%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08084'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  -2.02%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.746'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -0.61%  '

$ws.Range("E37").Value = '  -1.69%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2678'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  -2.76%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '10.70'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  -5.21%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09119'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  -1.48%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7602'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  -1.76%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.39'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  -3.33%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.420'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  -1.74%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.33'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  -0.90%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6966'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  -2.00%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.463'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -2.67%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.100'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -1.02%  '

$ws.Range("E48").Value = '  +0.14%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08327'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  -1.50%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '134.86'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  -1.01%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.220'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -3.62%  '
